# "added all place holder slides and updated spreadsheet"
#
# Spreadsheet-side change: a new "Order" helper column is inserted at the
# very left of the table (A), a priority/order number is typed in for each
# existing feature row, and the whole table is then sorted ascending by
# that new column (rows without a number - "Delta Optimization" - sort to
# the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A - shifts the existing A:E data to B:F.
[void]$ws.Columns.Item(1).Insert()

# Header for the new column, styled like the other bold header cells.
$ws.Range("A1").Value = "Order"
$ws.Range("A1").Font.Bold = $true

# Order numbers, one per existing data row (rows 2-13, in their current,
# pre-sort, top-to-bottom order). "Delta Optimization" (the row that ends
# up at r=7, i.e. index 5 below) never gets a number.
$orderValues = @(6, 3, 4, 5, 7, $null, 10, 2, 1, 9, 8, 7)

for ($i = 0; $i -lt $orderValues.Length; $i++) {
    $rowNum = 2 + $i
    $val = $orderValues[$i]
    if ($null -ne $val) {
        $ws.Cells.Item($rowNum, 1).Value = $val
    }
}

# Sort the table (including header) ascending by the new Order column.
$sortObj = $ws.Sort()
$sortFields = $sortObj.SortFields()
$sortFields.Clear()
[void]$sortFields.Add($ws.Range("A1:A13"))
$sortObj.SetRange($ws.Range("A1:F13"))
$sortObj.Header = 1
$sortObj.Apply()

# Match the author's final selection.
[void]$ws.Range("B12").Select()
